$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "PREDIOS INUNDACION"
$ws.Range("D3").Value = "PREDIOS"
$ws.Range("D4").Value = "ZONAS GEOECONOMICAS_U"
$ws.Range("D5").Value = "ZONAS GEOECONOMICAS-R"
$ws.Range("D6").Value = "ZONAS_FISICAS_U"
$ws.Range("D7").Value = "ZONAS_FISICAS_R"

$ws.Range("D8").Select()
